$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 21 (shifts old rows 21-75 down to 22-76,
# so the previously-last row 75 now also lives on as the new row 76).
$ws.Rows(21).Insert()

# Populate the newly inserted row 21 with this week's data. The
# categorical columns (Calidad/Unidad/Origen/etc.) mirror what used to be
# in the old row 21, only the date and price figures are new.
$ws.Range("A21").Value = 7
$ws.Range("B21").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C21").Value = "Ñuble"
$ws.Range("D21").Value = 44414
$ws.Range("E21").Value = 16
$ws.Range("F21").Value = 100112045
$ws.Range("G21").Value = "Zapallo"
$ws.Range("H21").Value = "Camote"
$ws.Range("I21").Value = "1a (guarda)"
$ws.Range("J21").Value = 300
$ws.Range("K21").Value = 400
$ws.Range("L21").Value = 450
$ws.Range("M21").Value = 425
$ws.Range("N21").Value = "$/kilo (volumen en unidades)"
$ws.Range("O21").Value = "Región del Maule"
$ws.Range("P21").Value = 425
$ws.Range("Q21").Value = 1
$ws.Range("R21").Value = "Hortaliza"
